# Add 6 new "false_S" pathway-gene rows to the "all_pathway_genes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_pathway_genes")

# Columns: A=_id  B=chel_query  C=chel_pathway  D=chel_gene_group
#          E=chel_gene  F=chel_subunit  G=_min_len  H=_max_len  I=_descr
$rows = @(
    @("Rv2208",          "-", "false_S", "-", "-", "adenosylcobinamide-GDP ribazoletransferase (Mycobacterium tuberculosis H37Rv)"),
    @("PP_1681",         "-", "false_S", "-", "-", "adenosylcobinamide-GDP ribazoletransferase (Pseudomonas putida KT2440)"),
    @("RD1_RS07205",     "-", "false_S", "-", "-", "CbbQ/NirQ/NorQ/GpvN family protein (Roseobacter denitrificans OCh 114)"),
    @("APZ15_RS35280",   "-", "false_S", "-", "-", "ATPase AAA (Burkholderia cepacia ATCC 25416)"),
    @("PNI01S_RS24580",  "-", "false_S", "-", "-", "CbbQ/NirQ/NorQ/GpvN family protein (Pseudomonas nitroreducens NBRC 12694)"),
    @("MAQ5080_RS04760", "-", "false_S", "-", "-", "ATPase (Marinomonas aquimarina)")
)

$startRow = 68

# First pass: populate _id / chel_query (A,B), chel_pathway (C), chel_subunit (F),
# _min_len/_max_len (G,H) and _descr (I) for every new row.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $id = $rows[$i][0]
    $pathway = $rows[$i][1]
    $subunit = $rows[$i][3]
    $lens = $rows[$i][4]
    $descr = $rows[$i][5]

    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $id
    $ws.Cells.Item($r, 3).Value = $pathway
    $ws.Cells.Item($r, 6).Value = $subunit
    $ws.Cells.Item($r, 7).Value = $lens
    $ws.Cells.Item($r, 8).Value = $lens
    $ws.Cells.Item($r, 9).Value = $descr
}

# Second pass: populate chel_gene_group / chel_gene (D,E) with "false_S" last,
# so the shared string is appended at the very end of the table.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $geneGroup = $rows[$i][2]

    $ws.Cells.Item($r, 4).Value = $geneGroup
    $ws.Cells.Item($r, 5).Value = $geneGroup
}

# Mirror the final cursor position left behind by the author (bottom of the
# newly appended block, scrolled so row 36 is the first visible row below
# the frozen header).
$lastRow = $startRow + $rows.Count - 1
$excel.ActiveWindow.ScrollRow = 36
$ws.Cells.Item($lastRow, 1).Select()
